$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = 95592
$ws.Range("B2").Value = "Dr. Thomas Moraes"
$ws.Range("C2").Value = "Jurídico"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("G2").Value = 9931.51

# Row 3
$ws.Range("A3").Value = 80078
$ws.Range("B3").Value = "André Cardoso"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Viagem de negócios"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45091
$ws.Range("G3").Value = 5452.42

# Row 4
$ws.Range("A4").Value = 48637
$ws.Range("B4").Value = "Lorena Carvalho"
$ws.Range("C4").Value = "Jurídico"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 45098
$ws.Range("G4").Value = 11672.37

# Row 5
$ws.Range("A5").Value = 61263
$ws.Range("B5").Value = "Sra. Emanuella Fogaça"
$ws.Range("C5").Value = "Marketing"
$ws.Range("D5").Value = "Doença"
$ws.Range("F5").Value = 45106
$ws.Range("G5").Value = 4956.17

# Row 6
$ws.Range("A6").Value = 76968
$ws.Range("B6").Value = "Maria Alice Rocha"
$ws.Range("C6").Value = "Operações"
$ws.Range("D6").Value = "Consulta médica"
$ws.Range("F6").Value = 45104
$ws.Range("G6").Value = 12476.3

# Row 7
$ws.Range("A7").Value = 43312
$ws.Range("B7").Value = "Caroline da Cruz"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45097
$ws.Range("G7").Value = 7462.74

# Row 8
$ws.Range("A8").Value = 40365
$ws.Range("B8").Value = "Marcela da Cruz"
$ws.Range("C8").Value = "Vendas"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45098
$ws.Range("G8").Value = 11973.47

# Row 9
$ws.Range("A9").Value = 20985
$ws.Range("B9").Value = "Elisa Gomes"
$ws.Range("C9").Value = "Jurídico"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45082
$ws.Range("G9").Value = 11066.15

# Row 10
$ws.Range("A10").Value = 24849
$ws.Range("B10").Value = "Heloísa Santos"
$ws.Range("C10").Value = "Jurídico"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45085
$ws.Range("G10").Value = 5944.88

# Row 11
$ws.Range("A11").Value = 19740
$ws.Range("B11").Value = "Davi Monteiro"
$ws.Range("C11").Value = "Operações"
$ws.Range("D11").Value = "Outros"
$ws.Range("F11").Value = 45099
$ws.Range("G11").Value = 8637.129999999999
